$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (date / text / wrap-text styles) from the last existing
# changelog row (16) down into the new row (17) before writing values, so
# the new cells pick up the same cell styles (s="1"/"2"/"3").
$ws.Range("A16:C16").Copy() | Out-Null
$ws.Range("A17:C17").PasteSpecial(-4122) | Out-Null

# New changelog entry: 1.4.1 (2021-06-06)
$ws.Range("A17").Value = 44353
$ws.Range("B17").Value = "1.4.1"
$ws.Range("C17").Value = "Bugs:`n- Visualization of coins graphic was not correct (not sorted along date)`n- Tracking of all burned DFI (manual, fees and unused block rewards)"

# Row height follows the same "15pt per wrapped line" convention as the
# other multi-line changelog rows (3 lines -> 45pt).
$ws.Rows.Item(17).RowHeight = 45

$excel.CutCopyMode = $false

# Move the active selection down past the newly added row, as it would be
# after typing the new entry.
$ws.Range("C18").Select() | Out-Null
